$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.756.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").Value = "'1.775.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.48%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = "'327.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = "'0.4575"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.29%  '
$ws.Range("D8").Value = "'0.3586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D9").Value = "'0.07489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").Value = "'41.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("D11").Value = "'1.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").Value = "'20.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").Value = "'6.049"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("D15").Value = "'7.222"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("D16").Value = "'1.775.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("D17").Value = "'93.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("D19").Value = "'0.06438"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.82%  '
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("E21").Value = '  +2.06%  '
$ws.Range("D22").Value = "'5.812"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.79%  '
$ws.Range("D23").Value = "'27.795.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = '  +1.18%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = "'164.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.51%  '
$ws.Range("E27").Value = '  -1.14%  '
$ws.Range("D28").Value = "'1.980.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.55%  '
$ws.Range("D29").Value = "'2.194"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.05%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").Value = "'1.104"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.25%  '
$ws.Range("D32").Value = "'0.09210"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.13%  '
$ws.Range("D33").Value = "'3.670"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.35%  '
$ws.Range("D34").Value = "'5.543"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").Value = "'11.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.49%  '
$ws.Range("D36").Value = "'0.02296"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").Value = "'0.06178"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.73%  '
$ws.Range("D38").Value = "'0.2090"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.26%  '
$ws.Range("D39").Value = "'0.6330"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").Value = "'4.967"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.47%  '
$ws.Range("D41").Value = "'1.186"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.67%  '
$ws.Range("D42").Value = "'1.388"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").Value = "'7.846"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.16%  '
$ws.Range("D44").Value = "'13.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("D47").Value = "'122.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("D48").Value = "'1.954"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").Value = "'0.06931"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("D51").Value = "'72.56"
$ws.Range("D51").Style = "Normal"
